$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '58.166.02'
$ws.Range('E2').Value = '  -3.74%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.450.61'
$ws.Range('E3').Value = '  -3.67%  '

$ws.Range('E4').Value = '  -0.01%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '527.97'
$ws.Range('E5').Value = '  -2.71%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '133.19'
$ws.Range('E6').Value = '  -8.83%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  +0.65%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.553'
$ws.Range('E8').Value = '  -3.14%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.452.98'
$ws.Range('E9').Value = '  -4.65%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0986'
$ws.Range('E10').Value = '  -3.25%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.160'
$ws.Range('E11').Value = '  -0.37%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.35'
$ws.Range('E12').Value = '  -3.49%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.343'
$ws.Range('E13').Value = '  -6.02%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.890.16'
$ws.Range('E14').Value = '  -3.82%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '58.129.17'
$ws.Range('E15').Value = '  -3.68%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '22.36'
$ws.Range('E16').Value = '  -8.50%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000137'
$ws.Range('E17').Value = '  -4.72%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.460.95'
$ws.Range('E18').Value = '  -5.13%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.61'
$ws.Range('E19').Value = '  -6.66%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.16'
$ws.Range('E20').Value = '  -4.88%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '317.51'
$ws.Range('E21').Value = '  -3.59%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.996'
$ws.Range('E22').Value = '  -0.17%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.66'
$ws.Range('E23').Value = '  -4.97%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '61.98'
$ws.Range('E24').Value = '  -1.62%  '

$ws.Range('B25').Value = 'Kaspa'
$ws.Range('C25').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.163'
$ws.Range('E25').Value = '  -2.71%  '

$ws.Range('B26').Value = 'Polygon'
$ws.Range('C26').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.402'
$ws.Range('E26').Value = '  -8.90%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.986'
$ws.Range('E27').Value = '  -0.98%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.49'
$ws.Range('E28').Value = '  -7.52%  '

$ws.Range('B29').Value = 'Aptos'
$ws.Range('C29').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.50'
$ws.Range('E29').Value = '  -9.26%  '

$ws.Range('B30').Value = 'PEPE'
$ws.Range('C30').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0₃0747'
$ws.Range('E30').Value = '  -7.08%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.75'
$ws.Range('E31').Value = '  -4.20%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '163.26'
$ws.Range('E32').Value = '  -0.21%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.999'
$ws.Range('E33').Value = '  +0.15%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.06'
$ws.Range('E34').Value = '  -12.28%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.34'
$ws.Range('E35').Value = '  -9.83%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '18.08'
$ws.Range('E36').Value = '  -4.05%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.99'
$ws.Range('E37').Value = '  -11.60%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.53'
$ws.Range('E38').Value = '  -7.52%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '36.36'
$ws.Range('E39').Value = '  -2.22%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.52'
$ws.Range('E40').Value = '  -6.48%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.776'
$ws.Range('E41').Value = '  -8.08%  '

$ws.Range('B42').Value = 'FirstDigitalUSD'
$ws.Range('C42').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.999'
$ws.Range('E42').Value = '  +0.87%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '271.21'
$ws.Range('E43').Value = '  -11.54%  '

$ws.Range('B44').Value = 'RenderToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '4.98'
$ws.Range('E44').Value = '  -13.08%  '

$ws.Range('E45').Value = '  -0.04%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.582'
$ws.Range('E46').Value = '  -4.86%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0923'
$ws.Range('E47').Value = '  -1.84%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '120.79'
$ws.Range('E48').Value = '  -4.79%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0500'
$ws.Range('E49').Value = '  -5.05%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0217'
$ws.Range('E50').Value = '  -5.94%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '16.86'
$ws.Range('E51').Value = '  -8.45%  '
